$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 2500  # H48: 4998 -> 2500
$ws.Cells.Item(48, 9).Value = 1000  # I48: 0 -> 1000
$ws.Cells.Item(48, 10).Value = 4000  # J48: 4998 -> 4000
$ws.Cells.Item(48, 11).Value = 3000  # K48: 0 -> 3000
$ws.Cells.Item(48, 12).Value = 12000  # L48: 14994 -> 12000
$ws.Cells.Item(48, 13).Value = -2708  # M48: None -> -2708
$ws.Cells.Item(48, 14).Value = -12584  # N48: -15578 -> -12584
$ws.Cells.Item(56, 8).Value = 2500  # H56: 4998 -> 2500
$ws.Cells.Item(56, 9).Value = 1000  # I56: 0 -> 1000
$ws.Cells.Item(56, 10).Value = 4000  # J56: 4998 -> 4000
$ws.Cells.Item(56, 11).Value = 3000  # K56: 0 -> 3000
$ws.Cells.Item(56, 12).Value = 12000  # L56: 14994 -> 12000
$ws.Cells.Item(56, 13).Value = -2466  # M56: None -> -2466
$ws.Cells.Item(56, 14).Value = -13068  # N56: -16062 -> -13068
$ws.Cells.Item(62, 8).Value = 7221  # H62: 7777 -> 7221
$ws.Cells.Item(62, 9).Value = 4997.5  # I62: 4998 -> 4997.5
$ws.Cells.Item(62, 11).Value = 4997.5  # K62: 4998 -> 4997.5
$ws.Cells.Item(62, 13).Value = -4373.5  # M62: -4374 -> -4373.5
$ws.Cells.Item(64, 8).Value = 4820  # H64: 4516.6665 -> 4820
$ws.Cells.Item(64, 10).Value = 5033.3335  # J64: 4525 -> 5033.3335
$ws.Cells.Item(64, 12).Value = 5033.3335  # L64: 4525 -> 5033.3335
$ws.Cells.Item(64, 14).Value = -5529.3335  # N64: -5021 -> -5529.3335
$ws.Cells.Item(65, 8).Value = 7221  # H65: 7777 -> 7221
$ws.Cells.Item(65, 9).Value = 4997.5  # I65: 4998 -> 4997.5
$ws.Cells.Item(65, 11).Value = 24987.5  # K65: 24990 -> 24987.5
$ws.Cells.Item(65, 13).Value = -21867.5  # M65: -21870 -> -21867.5
$ws.Cells.Item(67, 8).Value = 4820  # H67: 4516.6665 -> 4820
$ws.Cells.Item(67, 10).Value = 5033.3335  # J67: 4525 -> 5033.3335
$ws.Cells.Item(67, 12).Value = 5033.3335  # L67: 4525 -> 5033.3335
$ws.Cells.Item(67, 14).Value = -6749.3335  # N67: -6241 -> -6749.3335
$ws.Cells.Item(98, 8).Value = 2681.1333  # H98: 2684.2 -> 2681.1333
$ws.Cells.Item(98, 9).Value = 2446.4546  # I98: 2363 -> 2446.4546
$ws.Cells.Item(98, 10).Value = 3326.5  # J98: 3969 -> 3326.5
$ws.Cells.Item(98, 11).Value = 2446.4546  # K98: 2363 -> 2446.4546
$ws.Cells.Item(98, 12).Value = 3326.5  # L98: 3969 -> 3326.5
$ws.Cells.Item(98, 13).Value = -948.4546  # M98: -865 -> -948.4546
$ws.Cells.Item(98, 14).Value = -6322.5  # N98: -6965 -> -6322.5
$ws.Cells.Item(111, 8).Value = 50381.875  # H111: 54722.43 -> 50381.875
$ws.Cells.Item(111, 9).Value = 44007.855  # I111: 48009.5 -> 44007.855
$ws.Cells.Item(111, 11).Value = 132023.565  # K111: 144028.5 -> 132023.565
$ws.Cells.Item(111, 13).Value = -128956.565  # M111: -140961.5 -> -128956.565
$ws.Cells.Item(122, 8).Value = 2681.1333  # H122: 2684.2 -> 2681.1333
$ws.Cells.Item(122, 9).Value = 2446.4546  # I122: 2363 -> 2446.4546
$ws.Cells.Item(122, 10).Value = 3326.5  # J122: 3969 -> 3326.5
$ws.Cells.Item(122, 11).Value = 7339.3638  # K122: 7089 -> 7339.3638
$ws.Cells.Item(122, 12).Value = 9979.5  # L122: 11907 -> 9979.5
$ws.Cells.Item(122, 13).Value = -4889.3638  # M122: -4639 -> -4889.3638
$ws.Cells.Item(122, 14).Value = -14879.5  # N122: -16807 -> -14879.5
$ws.Cells.Item(137, 8).Value = 3120.1765  # H137: 3045.257 -> 3120.1765
$ws.Cells.Item(137, 9).Value = 2849.4285  # I137: 2768.3447 -> 2849.4285
$ws.Cells.Item(137, 11).Value = 8548.2855  # K137: 8305.034100000001 -> 8548.2855
$ws.Cells.Item(137, 13).Value = -5998.2855  # M137: -5755.034100000001 -> -5998.2855

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1421.75  # H2: 1456.8148 -> 1421.75
$ws.Cells.Item(2, 9).Value = 1451.2  # I2: 1491.875 -> 1451.2
$ws.Cells.Item(2, 11).Value = 1451.2  # K2: 1491.875 -> 1451.2
$ws.Cells.Item(2, 13).Value = -1338.2  # M2: -1378.875 -> -1338.2
$ws.Cells.Item(32, 8).Value = 1712.84  # H32: 1544.76 -> 1712.84
$ws.Cells.Item(32, 9).Value = 1712.84  # I32: 1544.76 -> 1712.84
$ws.Cells.Item(32, 11).Value = 1712.84  # K32: 1544.76 -> 1712.84
$ws.Cells.Item(32, 13).Value = -1425.84  # M32: -1257.76 -> -1425.84
$ws.Cells.Item(61, 8).Value = 6916.1885  # H61: 7381.735 -> 6916.1885
$ws.Cells.Item(61, 9).Value = 7076.9585  # I61: 7483.089 -> 7076.9585
$ws.Cells.Item(61, 10).Value = 5372.8  # J61: 6241.5 -> 5372.8
$ws.Cells.Item(61, 11).Value = 7076.9585  # K61: 7483.089 -> 7076.9585
$ws.Cells.Item(61, 12).Value = 5372.8  # L61: 6241.5 -> 5372.8
$ws.Cells.Item(61, 13).Value = -6864.9585  # M61: -7271.089 -> -6864.9585
$ws.Cells.Item(61, 14).Value = -5796.8  # N61: -6665.5 -> -5796.8
$ws.Cells.Item(74, 8).Value = 2725.3171  # H74: 2912.25 -> 2725.3171
$ws.Cells.Item(74, 9).Value = 1682.3462  # I74: 1729.2727 -> 1682.3462
$ws.Cells.Item(74, 10).Value = 4533.1333  # J74: 4771.2144 -> 4533.1333
$ws.Cells.Item(74, 11).Value = 1682.3462  # K74: 1729.2727 -> 1682.3462
$ws.Cells.Item(74, 12).Value = 4533.1333  # L74: 4771.2144 -> 4533.1333
$ws.Cells.Item(74, 13).Value = -808.3462  # M74: -855.2727 -> -808.3462
$ws.Cells.Item(74, 14).Value = -6281.1333  # N74: -6519.2144 -> -6281.1333
$ws.Cells.Item(77, 8).Value = 2725.3171  # H77: 2912.25 -> 2725.3171
$ws.Cells.Item(77, 9).Value = 1682.3462  # I77: 1729.2727 -> 1682.3462
$ws.Cells.Item(77, 10).Value = 4533.1333  # J77: 4771.2144 -> 4533.1333
$ws.Cells.Item(77, 11).Value = 8411.731  # K77: 8646.363499999999 -> 8411.731
$ws.Cells.Item(77, 12).Value = 22665.6665  # L77: 23856.072 -> 22665.6665
$ws.Cells.Item(77, 13).Value = -4043.731  # M77: -4278.363499999999 -> -4043.731
$ws.Cells.Item(77, 14).Value = -31401.6665  # N77: -32592.072 -> -31401.6665
$ws.Cells.Item(110, 8).Value = 1431.0769  # H110: 1483.75 -> 1431.0769
$ws.Cells.Item(110, 9).Value = 1270.6  # I110: 1323 -> 1270.6
$ws.Cells.Item(110, 11).Value = 1270.6  # K110: 1323 -> 1270.6
$ws.Cells.Item(110, 13).Value = 774.4000000000001  # M110: 722 -> 774.4000000000001
$ws.Cells.Item(116, 8).Value = 1421.75  # H116: 1456.8148 -> 1421.75
$ws.Cells.Item(116, 9).Value = 1451.2  # I116: 1491.875 -> 1451.2
$ws.Cells.Item(116, 11).Value = 1451.2  # K116: 1491.875 -> 1451.2
$ws.Cells.Item(116, 13).Value = 842.8  # M116: 802.125 -> 842.8
$ws.Cells.Item(132, 8).Value = 2164.9524  # H132: 2172.492 -> 2164.9524
$ws.Cells.Item(132, 9).Value = 1941.7234  # I132: 1951.8298 -> 1941.7234
$ws.Cells.Item(132, 11).Value = 5825.1702  # K132: 5855.4894 -> 5825.1702
$ws.Cells.Item(132, 13).Value = -3295.1702  # M132: -3325.4894 -> -3295.1702
$ws.Cells.Item(136, 8).Value = 6916.1885  # H136: 7381.735 -> 6916.1885
$ws.Cells.Item(136, 9).Value = 7076.9585  # I136: 7483.089 -> 7076.9585
$ws.Cells.Item(136, 10).Value = 5372.8  # J136: 6241.5 -> 5372.8
$ws.Cells.Item(136, 11).Value = 21230.8755  # K136: 22449.267 -> 21230.8755
$ws.Cells.Item(136, 12).Value = 16118.4  # L136: 18724.5 -> 16118.4
$ws.Cells.Item(136, 13).Value = -18680.8755  # M136: -19899.267 -> -18680.8755
$ws.Cells.Item(136, 14).Value = -21218.4  # N136: -23824.5 -> -21218.4

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1421.75  # H3: 1456.8148 -> 1421.75
$ws.Cells.Item(3, 9).Value = 1451.2  # I3: 1491.875 -> 1451.2
$ws.Cells.Item(3, 11).Value = 1451.2  # K3: 1491.875 -> 1451.2
$ws.Cells.Item(3, 13).Value = -1337.2  # M3: -1377.875 -> -1337.2
$ws.Cells.Item(64, 8).Value = 1266.875  # H64: 1266.9375 -> 1266.875
$ws.Cells.Item(64, 9).Value = 488  # I64: 488.5 -> 488
$ws.Cells.Item(64, 11).Value = 488  # K64: 488.5 -> 488
$ws.Cells.Item(64, 13).Value = -263  # M64: -263.5 -> -263
$ws.Cells.Item(67, 8).Value = 1266.875  # H67: 1266.9375 -> 1266.875
$ws.Cells.Item(67, 9).Value = 488  # I67: 488.5 -> 488
$ws.Cells.Item(67, 11).Value = 488  # K67: 488.5 -> 488
$ws.Cells.Item(67, 13).Value = 292  # M67: 291.5 -> 292
$ws.Cells.Item(76, 8).Value = 30000  # H76: 0 -> 30000
$ws.Cells.Item(76, 10).Value = 30000  # J76: 0 -> 30000
$ws.Cells.Item(76, 12).Value = 30000  # L76: 0 -> 30000
$ws.Cells.Item(76, 14).Value = -30630  # N76: None -> -30630
$ws.Cells.Item(79, 8).Value = 30000  # H79: 0 -> 30000
$ws.Cells.Item(79, 10).Value = 30000  # J79: 0 -> 30000
$ws.Cells.Item(79, 12).Value = 30000  # L79: 0 -> 30000
$ws.Cells.Item(79, 14).Value = -32184  # N79: None -> -32184
$ws.Cells.Item(105, 8).Value = 4523.577  # H105: 4747.375 -> 4523.577
$ws.Cells.Item(105, 9).Value = 3589.4119  # I105: 3822.9333 -> 3589.4119
$ws.Cells.Item(105, 11).Value = 3589.4119  # K105: 3822.9333 -> 3589.4119
$ws.Cells.Item(105, 13).Value = -1842.4119  # M105: -2075.9333 -> -1842.4119
$ws.Cells.Item(134, 8).Value = 3822.8333  # H134: 3913.3428 -> 3822.8333
$ws.Cells.Item(134, 9).Value = 3566.3333  # I134: 3657.3125 -> 3566.3333
$ws.Cells.Item(134, 11).Value = 10698.9999  # K134: 10971.9375 -> 10698.9999
$ws.Cells.Item(134, 13).Value = -8163.999899999999  # M134: -8436.9375 -> -8163.999899999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3929.5  # H62: 3814.2 -> 3929.5
$ws.Cells.Item(62, 10).Value = 4079.2222  # J62: 3891.3 -> 4079.2222
$ws.Cells.Item(62, 12).Value = 4079.2222  # L62: 3891.3 -> 4079.2222
$ws.Cells.Item(62, 14).Value = -5327.2222  # N62: -5139.3 -> -5327.2222
$ws.Cells.Item(65, 8).Value = 3929.5  # H65: 3814.2 -> 3929.5
$ws.Cells.Item(65, 10).Value = 4079.2222  # J65: 3891.3 -> 4079.2222
$ws.Cells.Item(65, 12).Value = 20396.111  # L65: 19456.5 -> 20396.111
$ws.Cells.Item(65, 14).Value = -26636.111  # N65: -25696.5 -> -26636.111
$ws.Cells.Item(99, 8).Value = 5876.9165  # H99: 5815 -> 5876.9165
$ws.Cells.Item(99, 9).Value = 5674.75  # I99: 5595.5884 -> 5674.75
$ws.Cells.Item(99, 11).Value = 5674.75  # K99: 5595.5884 -> 5674.75
$ws.Cells.Item(99, 13).Value = -4176.75  # M99: -4097.5884 -> -4176.75
$ws.Cells.Item(122, 8).Value = 4503.5293  # H122: 4541.7646 -> 4503.5293
$ws.Cells.Item(122, 9).Value = 4597.5  # I122: 4638.125 -> 4597.5
$ws.Cells.Item(122, 11).Value = 13792.5  # K122: 13914.375 -> 13792.5
$ws.Cells.Item(122, 13).Value = -11342.5  # M122: -11464.375 -> -11342.5
$ws.Cells.Item(126, 8).Value = 5876.9165  # H126: 5815 -> 5876.9165
$ws.Cells.Item(126, 9).Value = 5674.75  # I126: 5595.5884 -> 5674.75
$ws.Cells.Item(126, 11).Value = 17024.25  # K126: 16786.7652 -> 17024.25
$ws.Cells.Item(126, 13).Value = -14554.25  # M126: -14316.7652 -> -14554.25
$ws.Cells.Item(134, 8).Value = 2494.1702  # H134: 2590.4 -> 2494.1702
$ws.Cells.Item(134, 9).Value = 886.4857  # I134: 920.2727 -> 886.4857
$ws.Cells.Item(134, 11).Value = 2659.4571  # K134: 2760.8181 -> 2659.4571
$ws.Cells.Item(134, 13).Value = -124.4570999999996  # M134: -225.8181 -> -124.4570999999996
$ws.Cells.Item(141, 8).Value = 547820.7  # H141: 688562.9 -> 547820.7
$ws.Cells.Item(141, 10).Value = 547820.7  # J141: 688562.9 -> 547820.7
$ws.Cells.Item(141, 12).Value = 547820.7  # L141: 688562.9 -> 547820.7
$ws.Cells.Item(141, 14).Value = -558180.7  # N141: -698922.9 -> -558180.7

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 464  # H20: 464.33334 -> 464
$ws.Cells.Item(20, 9).Value = 464  # I20: 464.33334 -> 464
$ws.Cells.Item(20, 11).Value = 1392  # K20: 1393.00002 -> 1392
$ws.Cells.Item(20, 13).Value = -1165  # M20: -1166.00002 -> -1165
$ws.Cells.Item(23, 8).Value = 166668420  # H23: 200001920 -> 166668420
$ws.Cells.Item(23, 9).Value = 162  # I23: 163 -> 162
$ws.Cells.Item(23, 10).Value = 200002060  # J23: 250002350 -> 200002060
$ws.Cells.Item(23, 11).Value = 486  # K23: 489 -> 486
$ws.Cells.Item(23, 12).Value = 600006180  # L23: 750007050 -> 600006180
$ws.Cells.Item(23, 13).Value = -251  # M23: -254 -> -251
$ws.Cells.Item(23, 14).Value = -600006650  # N23: -750007520 -> -600006650
$ws.Cells.Item(25, 8).Value = 1337.2  # H25: 1337.4 -> 1337.2
$ws.Cells.Item(25, 9).Value = 295.33334  # I25: 295.66666 -> 295.33334
$ws.Cells.Item(25, 11).Value = 886.0000200000001  # K25: 886.9999799999999 -> 886.0000200000001
$ws.Cells.Item(25, 13).Value = -717.0000200000001  # M25: -717.9999799999999 -> -717.0000200000001
$ws.Cells.Item(30, 8).Value = 1337.2  # H30: 1337.4 -> 1337.2
$ws.Cells.Item(30, 9).Value = 295.33334  # I30: 295.66666 -> 295.33334
$ws.Cells.Item(30, 11).Value = 886.0000200000001  # K30: 886.9999799999999 -> 886.0000200000001
$ws.Cells.Item(30, 13).Value = -784.0000200000001  # M30: -784.9999799999999 -> -784.0000200000001
$ws.Cells.Item(33, 8).Value = 2023.3572  # H33: 2173.5386 -> 2023.3572
$ws.Cells.Item(33, 9).Value = 259.6  # I33: 260 -> 259.6
$ws.Cells.Item(33, 10).Value = 3003.2222  # J33: 3369.5 -> 3003.2222
$ws.Cells.Item(33, 11).Value = 1557.6  # K33: 1560 -> 1557.6
$ws.Cells.Item(33, 12).Value = 18019.3332  # L33: 20217 -> 18019.3332
$ws.Cells.Item(33, 13).Value = -1274.6  # M33: -1277 -> -1274.6
$ws.Cells.Item(33, 14).Value = -18585.3332  # N33: -20783 -> -18585.3332
$ws.Cells.Item(116, 8).Value = 6000  # H116: 3666 -> 6000
$ws.Cells.Item(116, 9).Value = 6000  # I116: 2500 -> 6000
$ws.Cells.Item(116, 10).Value = 6000  # J116: 5998 -> 6000
$ws.Cells.Item(116, 11).Value = 18000  # K116: 7500 -> 18000
$ws.Cells.Item(116, 12).Value = 18000  # L116: 17994 -> 18000
$ws.Cells.Item(116, 13).Value = -14558  # M116: -4058 -> -14558
$ws.Cells.Item(116, 14).Value = -24884  # N116: -24878 -> -24884

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0  # H5: 20000 -> 0
$ws.Cells.Item(5, 10).Value = 0  # J5: 20000 -> 0
$ws.Cells.Item(5, 12).Value = 0  # L5: 20000 -> 0
$ws.Cells.Item(5, 14).ClearContents()  # N5 (was -20224)
$ws.Cells.Item(80, 8).Value = 4614.727  # H80: 4614.8184 -> 4614.727
$ws.Cells.Item(80, 9).Value = 4583.4  # I80: 4583.6 -> 4583.4
$ws.Cells.Item(80, 11).Value = 4583.4  # K80: 4583.6 -> 4583.4
$ws.Cells.Item(80, 13).Value = -3585.4  # M80: -3585.6 -> -3585.4
$ws.Cells.Item(83, 8).Value = 4614.727  # H83: 4614.8184 -> 4614.727
$ws.Cells.Item(83, 9).Value = 4583.4  # I83: 4583.6 -> 4583.4
$ws.Cells.Item(83, 11).Value = 22917  # K83: 22918 -> 22917
$ws.Cells.Item(83, 13).Value = -17925  # M83: -17926 -> -17925
$ws.Cells.Item(102, 8).Value = 17922.719  # H102: 18478.484 -> 17922.719
$ws.Cells.Item(102, 9).Value = 2157.3333  # I102: 2213.6155 -> 2157.3333
$ws.Cells.Item(102, 11).Value = 2157.3333  # K102: 2213.6155 -> 2157.3333
$ws.Cells.Item(102, 13).Value = -535.3332999999998  # M102: -591.6154999999999 -> -535.3332999999998

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 17157.223  # H46: 17204.334 -> 17157.223
$ws.Cells.Item(46, 9).Value = 7331.6665  # I46: 10500 -> 7331.6665
$ws.Cells.Item(46, 10).Value = 22070  # J46: 19119.857 -> 22070
$ws.Cells.Item(46, 11).Value = 7331.6665  # K46: 10500 -> 7331.6665
$ws.Cells.Item(46, 12).Value = 22070  # L46: 19119.857 -> 22070
$ws.Cells.Item(46, 13).Value = -7143.6665  # M46: -10312 -> -7143.6665
$ws.Cells.Item(46, 14).Value = -22446  # N46: -19495.857 -> -22446
$ws.Cells.Item(95, 8).Value = 66991  # H95: 58593 -> 66991
$ws.Cells.Item(95, 10).Value = 64323  # J95: 54492.5 -> 64323
$ws.Cells.Item(95, 12).Value = 64323  # L95: 54492.5 -> 64323
$ws.Cells.Item(95, 14).Value = -69815  # N95: -59984.5 -> -69815

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 5000  # H8: 3000 -> 5000
$ws.Cells.Item(8, 9).Value = 0  # I8: 1000 -> 0
$ws.Cells.Item(8, 11).Value = 0  # K8: 1000 -> 0
$ws.Cells.Item(8, 13).ClearContents()  # M8 (was -860)
$ws.Cells.Item(96, 8).Value = 5274.75  # H96: 3983.1667 -> 5274.75
$ws.Cells.Item(96, 9).Value = 5549.5  # I96: 4133 -> 5549.5
$ws.Cells.Item(96, 10).Value = 5000  # J96: 3833.3333 -> 5000
$ws.Cells.Item(96, 11).Value = 5549.5  # K96: 4133 -> 5549.5
$ws.Cells.Item(96, 12).Value = 5000  # L96: 3833.3333 -> 5000
$ws.Cells.Item(96, 13).Value = -4176.5  # M96: -2760 -> -4176.5
$ws.Cells.Item(96, 14).Value = -7746  # N96: -6579.3333 -> -7746
$ws.Cells.Item(132, 8).Value = 1497.8  # H132: 1525.8276 -> 1497.8
$ws.Cells.Item(132, 9).Value = 1358.6666  # I132: 1384.5769 -> 1358.6666
$ws.Cells.Item(132, 11).Value = 4075.9998  # K132: 4153.7307 -> 4075.9998
$ws.Cells.Item(132, 13).Value = -1545.9998  # M132: -1623.7307 -> -1545.9998
$ws.Cells.Item(136, 8).Value = 8430.925999999999  # H136: 9436.333000000001 -> 8430.925999999999
$ws.Cells.Item(136, 9).Value = 9484.392  # I136: 10848.9 -> 9484.392
$ws.Cells.Item(136, 11).Value = 28453.176  # K136: 32546.7 -> 28453.176
$ws.Cells.Item(136, 13).Value = -25903.176  # M136: -29996.7 -> -25903.176
